# Auto-generated edit script: update Price (D) and Volume(1h) (E) columns
# to reflect refreshed crypto market data.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($ws, $cellRef, $text) {
    $cell = $ws.Range($cellRef)
    $cell.NumberFormat = "@"
    $cell.Value = $text
    $cell.Style = "Normal"
}

Set-TextValue $ws "D2" "23.513.43"
Set-TextValue $ws "E2" "  -0.29%  "
Set-TextValue $ws "D3" "1.639.78"
Set-TextValue $ws "E3" "  -0.42%  "
Set-TextValue $ws "D4" "0.9990"
Set-TextValue $ws "E4" "  -0.54%  "
Set-TextValue $ws "D5" "0.9996"
Set-TextValue $ws "E5" "  -0.22%  "
Set-TextValue $ws "D6" "304.26"
Set-TextValue $ws "E6" "  -0.38%  "
Set-TextValue $ws "D7" "0.3797"
Set-TextValue $ws "E7" "  +0.55%  "
Set-TextValue $ws "D8" "51.67"
Set-TextValue $ws "E8" "  -2.38%  "
Set-TextValue $ws "E9" "  -1.33%  "
Set-TextValue $ws "D10" "0.08194"
Set-TextValue $ws "E10" "  +0.44%  "
Set-TextValue $ws "E11" "  -2.51%  "
Set-TextValue $ws "D12" "1.000"
Set-TextValue $ws "E12" "  -0.39%  "
Set-TextValue $ws "D13" "22.51"
Set-TextValue $ws "E13" "  -2.69%  "
Set-TextValue $ws "D14" "6.472"
Set-TextValue $ws "E14" "  -3.59%  "
Set-TextValue $ws "D15" "7.404"
Set-TextValue $ws "E15" "  +0.15%  "
Set-TextValue $ws "E16" "  -2.27%  "
Set-TextValue $ws "D17" "1.633.16"
Set-TextValue $ws "E17" "  -1.02%  "
Set-TextValue $ws "D18" "95.43"
Set-TextValue $ws "E18" "  +0.32%  "
Set-TextValue $ws "D19" "0.06937"
Set-TextValue $ws "E19" "  +0.22%  "
Set-TextValue $ws "D20" "6.597"
Set-TextValue $ws "E20" "  -0.06%  "
Set-TextValue $ws "D21" "17.52"
Set-TextValue $ws "E21" "  -4.59%  "
Set-TextValue $ws "D22" "0.9990"
Set-TextValue $ws "E22" "  -0.36%  "
Set-TextValue $ws "D23" "12.54"
Set-TextValue $ws "E23" "  -3.33%  "
Set-TextValue $ws "D24" "23.521.46"
Set-TextValue $ws "E24" "  -0.34%  "
Set-TextValue $ws "D25" "2.495"
Set-TextValue $ws "D26" "3.067"
Set-TextValue $ws "E26" "  -5.79%  "
Set-TextValue $ws "D27" "21.14"
Set-TextValue $ws "E27" "  -1.73%  "
Set-TextValue $ws "D28" "151.53"
Set-TextValue $ws "E28" "  -0.38%  "
Set-TextValue $ws "D29" "5.269"
Set-TextValue $ws "E29" "  -1.05%  "
Set-TextValue $ws "D30" "133.26"
Set-TextValue $ws "E30" "  -3.11%  "
Set-TextValue $ws "D31" "1.816.07"
Set-TextValue $ws "E31" "  -0.94%  "
Set-TextValue $ws "D32" "2.189"
Set-TextValue $ws "E32" "  -5.27%  "
Set-TextValue $ws "D33" "6.667"
Set-TextValue $ws "E33" "  -5.23%  "
Set-TextValue $ws "D34" "1.061"
Set-TextValue $ws "E34" "  +8.40%  "
Set-TextValue $ws "D35" "11.39"
Set-TextValue $ws "E35" "  +3.51%  "
Set-TextValue $ws "D36" "0.02768"
Set-TextValue $ws "E36" "  -4.26%  "
Set-TextValue $ws "D37" "0.2502"
Set-TextValue $ws "E37" "  -3.30%  "
Set-TextValue $ws "D38" "0.08785"
Set-TextValue $ws "D39" "0.07123"
Set-TextValue $ws "E39" "  -2.90%  "
Set-TextValue $ws "D40" "6.032"
Set-TextValue $ws "E40" "  -5.35%  "
Set-TextValue $ws "D41" "0.7054"
Set-TextValue $ws "E41" "  -2.25%  "
Set-TextValue $ws "E42" "  -3.27%  "
Set-TextValue $ws "E43" "  -4.53%  "
Set-TextValue $ws "D44" "12.16"
Set-TextValue $ws "E44" "  -4.66%  "
Set-TextValue $ws "D45" "0.6544"
Set-TextValue $ws "E45" "  -1.62%  "
Set-TextValue $ws "D46" "0.9990"
Set-TextValue $ws "E46" "  -0.27%  "
Set-TextValue $ws "D47" "2.288"
Set-TextValue $ws "E47" "  -4.37%  "
Set-TextValue $ws "D48" "3.968"
Set-TextValue $ws "E48" "  -1.34%  "
Set-TextValue $ws "E49" "  -0.90%  "
Set-TextValue $ws "D50" "128.66"
Set-TextValue $ws "E50" "  -0.25%  "
Set-TextValue $ws "D51" "1.190"
Set-TextValue $ws "E51" "  -3.27%  "
